# Applies the edit described by the diff: the data rows 4, 5 and 6 on the
# "Artfynd" sheet get their content cyclically rotated:
#   new row 4  <=  old row 5
#   new row 5  <=  old row 6
#   new row 6  <=  old row 4
# (row numbers/positions stay the same, only the field values move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values participate in the rotation, in sheet order.
$cols = @("A","B","E","F","G","H","Q","R","Z","AB","AC","AX")

# Snapshot current values for rows 4, 5 and 6 before overwriting anything.
$snapshot = @{}
foreach ($row in 4,5,6) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Mapping: destination row -> source row
$rotation = @{ 4 = 5; 5 = 6; 6 = 4 }

foreach ($destRow in 4,5,6) {
    $srcRow = $rotation[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$destRow")
        $val = $srcVals[$col]
        if ($col -eq "AC" -and ($val -eq $null -or $val -eq "")) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

$wb.Save()
